# =====================================================================
# Update infrastructure-costs.xlsx to November 2025 category structure
# =====================================================================
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Cover sheet: bump the "Solution" date
# ---------------------------------------------------------------
$wsCover = $wb.Worksheets.Item("Cover")
$wsCover.Range("C7").Value = "November 25, 2025"

# ---------------------------------------------------------------
# 2) Sizing Guidelines: grow from 7 data rows (3-9) to 9 data rows (3-11)
# ---------------------------------------------------------------
$wsSizing = $wb.Worksheets.Item("Sizing Guidelines")
if ($wsSizing.AutoFilterMode) { $wsSizing.AutoFilterMode = $false }

$wsSizing.Rows.Item(10).Insert()
$wsSizing.Range("A8:F8").Copy($wsSizing.Range("A10:F10"))
$wsSizing.Rows.Item(10).RowHeight = 26

$wsSizing.Rows.Item(11).Insert()
$wsSizing.Range("A9:F9").Copy($wsSizing.Range("A11:F11"))
$wsSizing.Rows.Item(11).RowHeight = 26

$wsSizing.Range("A3").Value = "Hardware"
$wsSizing.Range("B3").Value = "Network Devices"
$wsSizing.Range("C3").Value = "0-2 devices"
$wsSizing.Range("D3").Value = "3-5 devices"
$wsSizing.Range("E3").Value = "10+ devices"
$wsSizing.Range("F3").Value = "Based on site count"

$wsSizing.Range("A4").Value = "Cloud Services"
$wsSizing.Range("B4").Value = "Compute Instances"
$wsSizing.Range("C4").Value = "2-4 instances"
$wsSizing.Range("D4").Value = "5-10 instances"
$wsSizing.Range("E4").Value = "15+ instances"
$wsSizing.Range("F4").Value = "Production workload capacity"

$wsSizing.Range("A5").Value = "Cloud Services"
$wsSizing.Range("B5").Value = "Storage Volume"
$wsSizing.Range("C5").Value = "100 GB"
$wsSizing.Range("D5").Value = "500 GB - 1 TB"
$wsSizing.Range("E5").Value = "2+ TB"
$wsSizing.Range("F5").Value = "Based on data requirements"

$wsSizing.Range("A6").Value = "Cloud Services"
$wsSizing.Range("B6").Value = "Database Size"
$wsSizing.Range("C6").Value = "Small/Standard tier"
$wsSizing.Range("D6").Value = "Medium tier"
$wsSizing.Range("E6").Value = "Large/Enterprise tier"
$wsSizing.Range("F6").Value = "Based on transaction volume"

$wsSizing.Range("A7").Value = "Software Licenses"
$wsSizing.Range("B7").Value = "User Licenses"
$wsSizing.Range("C7").Value = "10-50 users"
$wsSizing.Range("D7").Value = "50-250 users"
$wsSizing.Range("E7").Value = "500+ users"
$wsSizing.Range("F7").Value = "Concurrent or named users"

$wsSizing.Range("A8").Value = "Software Licenses"
$wsSizing.Range("B8").Value = "Monitoring Tools"
$wsSizing.Range("C8").Value = "Basic monitoring"
$wsSizing.Range("D8").Value = "Advanced APM"
$wsSizing.Range("E8").Value = "Enterprise observability"
$wsSizing.Range("F8").Value = "Application performance monitoring"

$wsSizing.Range("A9").Value = "Connectivity"
$wsSizing.Range("B9").Value = "Network Bandwidth"
$wsSizing.Range("C9").Value = "100 Mbps"
$wsSizing.Range("D9").Value = "500 Mbps"
$wsSizing.Range("E9").Value = "1+ Gbps"
$wsSizing.Range("F9").Value = "WAN/Internet circuit speed"

$wsSizing.Range("A10").Value = "Support & Maintenance"
$wsSizing.Range("B10").Value = "Support Level"
$wsSizing.Range("C10").Value = "Business hours"
$wsSizing.Range("D10").Value = "24x5 support"
$wsSizing.Range("E10").Value = "24x7 premium"
$wsSizing.Range("F10").Value = "SLA requirements"

$wsSizing.Range("A11").Value = "Facilities"
$wsSizing.Range("B11").Value = "Data Center Space"
$wsSizing.Range("C11").Value = "N/A"
$wsSizing.Range("D11").Value = "N/A"
$wsSizing.Range("E11").Value = "N/A"
$wsSizing.Range("F11").Value = "Cloud-only solution (no facilities)"

$wsSizing.Range("A2:F11").AutoFilter()

# ---------------------------------------------------------------
# 3) Infrastructure Costs: grow from 8 data rows + TOTAL (3-10, 11)
#    to 10 data rows + TOTAL (3-12, 13)
# ---------------------------------------------------------------
$wsInfra = $wb.Worksheets.Item("Infrastructure Costs")
if ($wsInfra.AutoFilterMode) { $wsInfra.AutoFilterMode = $false }

$wsInfra.Rows.Item(11).Insert()
$wsInfra.Range("A9:K9").Copy($wsInfra.Range("A11:K11"))
$wsInfra.Rows.Item(11).RowHeight = 26

$wsInfra.Rows.Item(12).Insert()
$wsInfra.Range("A10:K10").Copy($wsInfra.Range("A12:K12"))
$wsInfra.Rows.Item(12).RowHeight = 26

# Row 3: Hardware / Network Switch
$wsInfra.Range("A3").Value = "Hardware"
$wsInfra.Range("B3").Value = "Network Switch"
$wsInfra.Range("C3").Value = "Network equipment"
$wsInfra.Range("D3").Value = 0
$wsInfra.Range("E3").Value = "Device"
$wsInfra.Range("F3").Value = 2500
$wsInfra.Range("G3").Formula = "=D3*F3"
$wsInfra.Range("H3").Formula = "=0"
$wsInfra.Range("I3").Formula = "=0"
$wsInfra.Range("J3").Formula = "=G3+H3+I3"
$wsInfra.Range("K3").Value = "Optional on-premises"

# Row 4: Cloud Services / Compute Instance
$wsInfra.Range("A4").Value = "Cloud Services"
$wsInfra.Range("B4").Value = "Compute Instance"
$wsInfra.Range("C4").Value = "Cloud VM or equivalent"
$wsInfra.Range("D4").Value = 2
$wsInfra.Range("E4").Value = "Instance/Month"
$wsInfra.Range("F4").Value = 150
$wsInfra.Range("G4").Formula = "=D4*F4*12"
$wsInfra.Range("H4").Formula = "=G4"
$wsInfra.Range("I4").Formula = "=G4"
$wsInfra.Range("J4").Formula = "=G4+H4+I4"
$wsInfra.Range("K4").Value = "Production workloads"

# Row 5: Cloud Services / Database
$wsInfra.Range("A5").Value = "Cloud Services"
$wsInfra.Range("B5").Value = "Database"
$wsInfra.Range("C5").Value = "Managed database service"
$wsInfra.Range("D5").Value = 1
$wsInfra.Range("E5").Value = "Instance/Month"
$wsInfra.Range("F5").Value = 200
$wsInfra.Range("G5").Formula = "=D5*F5*12"
$wsInfra.Range("H5").Formula = "=G5"
$wsInfra.Range("I5").Formula = "=G5"
$wsInfra.Range("J5").Formula = "=G5+H5+I5"
$wsInfra.Range("K5").Value = "Relational database"

# Row 6: Cloud Services / Storage
$wsInfra.Range("A6").Value = "Cloud Services"
$wsInfra.Range("B6").Value = "Storage"
$wsInfra.Range("C6").Value = "Object storage"
$wsInfra.Range("D6").Value = 100
$wsInfra.Range("E6").Value = "GB/Month"
$wsInfra.Range("F6").Value = 0.023
$wsInfra.Range("G6").Formula = "=D6*F6*12"
$wsInfra.Range("H6").Formula = "=G6"
$wsInfra.Range("I6").Formula = "=G6"
$wsInfra.Range("J6").Formula = "=G6+H6+I6"
$wsInfra.Range("K6").Value = "Data storage"

# Row 7: Cloud Services / Data Transfer
$wsInfra.Range("A7").Value = "Cloud Services"
$wsInfra.Range("B7").Value = "Data Transfer"
$wsInfra.Range("C7").Value = "Outbound transfer"
$wsInfra.Range("D7").Value = 50
$wsInfra.Range("E7").Value = "GB/Month"
$wsInfra.Range("F7").Value = 0.09
$wsInfra.Range("G7").Formula = "=D7*F7*12"
$wsInfra.Range("H7").Formula = "=G7"
$wsInfra.Range("I7").Formula = "=G7"
$wsInfra.Range("J7").Formula = "=G7+H7+I7"
$wsInfra.Range("K7").Value = "Egress charges"

# Row 8: Software Licenses / Monitoring
$wsInfra.Range("A8").Value = "Software Licenses"
$wsInfra.Range("B8").Value = "Monitoring"
$wsInfra.Range("C8").Value = "Application monitoring"
$wsInfra.Range("D8").Value = 1
$wsInfra.Range("E8").Value = "Service/Month"
$wsInfra.Range("F8").Value = 50
$wsInfra.Range("G8").Formula = "=D8*F8*12"
$wsInfra.Range("H8").Formula = "=G8"
$wsInfra.Range("I8").Formula = "=G8"
$wsInfra.Range("J8").Formula = "=G8+H8+I8"
$wsInfra.Range("K8").Value = "Metrics and logging"

# Row 9: Software Licenses / Security
$wsInfra.Range("A9").Value = "Software Licenses"
$wsInfra.Range("B9").Value = "Security"
$wsInfra.Range("C9").Value = "Web application firewall"
$wsInfra.Range("D9").Value = 1
$wsInfra.Range("E9").Value = "Service/Month"
$wsInfra.Range("F9").Value = 25
$wsInfra.Range("G9").Formula = "=D9*F9*12"
$wsInfra.Range("H9").Formula = "=G9"
$wsInfra.Range("I9").Formula = "=G9"
$wsInfra.Range("J9").Formula = "=G9+H9+I9"
$wsInfra.Range("K9").Value = "Security controls"

# Row 10: Connectivity / Internet Circuit
$wsInfra.Range("A10").Value = "Connectivity"
$wsInfra.Range("B10").Value = "Internet Circuit"
$wsInfra.Range("C10").Value = "Business internet"
$wsInfra.Range("D10").Value = 1
$wsInfra.Range("E10").Value = "Circuit/Month"
$wsInfra.Range("F10").Value = 500
$wsInfra.Range("G10").Formula = "=D10*F10*12"
$wsInfra.Range("H10").Formula = "=G10"
$wsInfra.Range("I10").Formula = "=G10"
$wsInfra.Range("J10").Formula = "=G10+H10+I10"
$wsInfra.Range("K10").Value = "Cloud-only (no circuit needed)"

# Row 11: Support & Maintenance / Cloud Support
$wsInfra.Range("A11").Value = "Support & Maintenance"
$wsInfra.Range("B11").Value = "Cloud Support"
$wsInfra.Range("C11").Value = "Provider support plan"
$wsInfra.Range("D11").Value = 1
$wsInfra.Range("E11").Value = "Plan/Month"
$wsInfra.Range("F11").Value = 100
$wsInfra.Range("G11").Formula = "=D11*F11*12"
$wsInfra.Range("H11").Formula = "=G11"
$wsInfra.Range("I11").Formula = "=G11"
$wsInfra.Range("J11").Formula = "=G11+H11+I11"
$wsInfra.Range("K11").Value = "Technical support"

# Row 12: Facilities / Data Center
$wsInfra.Range("A12").Value = "Facilities"
$wsInfra.Range("B12").Value = "Data Center"
$wsInfra.Range("C12").Value = "Rack space and power"
$wsInfra.Range("D12").Value = 0
$wsInfra.Range("E12").Value = "Rack/Month"
$wsInfra.Range("F12").Value = 0
$wsInfra.Range("G12").Formula = "=D12*F12*12"
$wsInfra.Range("H12").Formula = "=G12"
$wsInfra.Range("I12").Formula = "=G12"
$wsInfra.Range("J12").Formula = "=G12+H12+I12"
$wsInfra.Range("K12").Value = "Cloud-only (no facilities)"

# Row 13: TOTAL (now shifted down from row 11)
$wsInfra.Range("G13").Formula = "=SUM(G3:G12)"
$wsInfra.Range("H13").Formula = "=SUM(H3:H12)"
$wsInfra.Range("I13").Formula = "=SUM(I3:I12)"
$wsInfra.Range("J13").Formula = "=SUM(J3:J12)"

$wsInfra.Range("A2:K13").AutoFilter()

# ---------------------------------------------------------------
# 4) Credits: grow from 4 data rows (3-6) to 6 data rows (3-8)
# ---------------------------------------------------------------
$wsCredits = $wb.Worksheets.Item("Credits")
if ($wsCredits.AutoFilterMode) { $wsCredits.AutoFilterMode = $false }

$wsCredits.Rows.Item(7).Insert()
$wsCredits.Range("A5:D5").Copy($wsCredits.Range("A7:D7"))
$wsCredits.Rows.Item(7).RowHeight = 26

$wsCredits.Rows.Item(8).Insert()
$wsCredits.Range("A6:D6").Copy($wsCredits.Range("A8:D8"))
$wsCredits.Rows.Item(8).RowHeight = 26

$wsCredits.Range("A3").Value = "Hardware"
$wsCredits.Range("B3").Value = "Equipment Credit"
$wsCredits.Range("C3").Value = 0
$wsCredits.Range("D3").Value = "No hardware credits available"

$wsCredits.Range("A4").Value = "Cloud Services"
$wsCredits.Range("B4").Value = "Provider Credit"
$wsCredits.Range("C4").Value = -1800
$wsCredits.Range("D4").Value = "30% credit on eligible compute and database services"

$wsCredits.Range("A5").Value = "Software Licenses"
$wsCredits.Range("B5").Value = "Partner Credit"
$wsCredits.Range("C5").Value = 0
$wsCredits.Range("D5").Value = "No software credits available"

$wsCredits.Range("A6").Value = "Connectivity"
$wsCredits.Range("B6").Value = "Circuit Credit"
$wsCredits.Range("C6").Value = 0
$wsCredits.Range("D6").Value = "No connectivity credits available"

$wsCredits.Range("A7").Value = "Support & Maintenance"
$wsCredits.Range("B7").Value = "Program Credit"
$wsCredits.Range("C7").Value = 0
$wsCredits.Range("D7").Value = "No support credits available"

$wsCredits.Range("A8").Value = "Facilities"
$wsCredits.Range("B8").Value = "Facilities Credit"
$wsCredits.Range("C8").Value = 0
$wsCredits.Range("D8").Value = "No facilities (cloud-only solution)"

$wsCredits.Range("A2:D8").AutoFilter()

# ---------------------------------------------------------------
# 5) 3-Year Summary: grow from 4 data rows + TOTAL (3-6, 7)
#    to 6 data rows + TOTAL (3-8, 9)
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("3-Year Summary")
if ($wsSummary.AutoFilterMode) { $wsSummary.AutoFilterMode = $false }

$wsSummary.Rows.Item(7).Insert()
$wsSummary.Range("A5:G5").Copy($wsSummary.Range("A7:G7"))
$wsSummary.Rows.Item(7).RowHeight = 26

$wsSummary.Rows.Item(8).Insert()
$wsSummary.Range("A6:G6").Copy($wsSummary.Range("A8:G8"))
$wsSummary.Rows.Item(8).RowHeight = 26

$wsSummary.Range("A3").Value = "Hardware"
$wsSummary.Range("B3").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A3,'Infrastructure Costs'!`$G:`$G)"
$wsSummary.Range("C3").Formula = "=SUMIF(Credits!`$A:`$A,A3,Credits!`$C:`$C)"
$wsSummary.Range("D3").Formula = "=B3+C3"
$wsSummary.Range("E3").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A3,'Infrastructure Costs'!`$H:`$H)"
$wsSummary.Range("F3").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A3,'Infrastructure Costs'!`$I:`$I)"
$wsSummary.Range("G3").Formula = "=D3+E3+F3"

$wsSummary.Range("A4").Value = "Cloud Services"
$wsSummary.Range("B4").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A4,'Infrastructure Costs'!`$G:`$G)"
$wsSummary.Range("C4").Formula = "=SUMIF(Credits!`$A:`$A,A4,Credits!`$C:`$C)"
$wsSummary.Range("D4").Formula = "=B4+C4"
$wsSummary.Range("E4").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A4,'Infrastructure Costs'!`$H:`$H)"
$wsSummary.Range("F4").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A4,'Infrastructure Costs'!`$I:`$I)"
$wsSummary.Range("G4").Formula = "=D4+E4+F4"

$wsSummary.Range("A5").Value = "Software Licenses"
$wsSummary.Range("B5").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A5,'Infrastructure Costs'!`$G:`$G)"
$wsSummary.Range("C5").Formula = "=SUMIF(Credits!`$A:`$A,A5,Credits!`$C:`$C)"
$wsSummary.Range("D5").Formula = "=B5+C5"
$wsSummary.Range("E5").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A5,'Infrastructure Costs'!`$H:`$H)"
$wsSummary.Range("F5").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A5,'Infrastructure Costs'!`$I:`$I)"
$wsSummary.Range("G5").Formula = "=D5+E5+F5"

$wsSummary.Range("A6").Value = "Connectivity"
$wsSummary.Range("B6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$G:`$G)"
$wsSummary.Range("C6").Formula = "=SUMIF(Credits!`$A:`$A,A6,Credits!`$C:`$C)"
$wsSummary.Range("D6").Formula = "=B6+C6"
$wsSummary.Range("E6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$H:`$H)"
$wsSummary.Range("F6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$I:`$I)"
$wsSummary.Range("G6").Formula = "=D6+E6+F6"

$wsSummary.Range("A7").Value = "Support & Maintenance"
$wsSummary.Range("B7").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A7,'Infrastructure Costs'!`$G:`$G)"
$wsSummary.Range("C7").Formula = "=SUMIF(Credits!`$A:`$A,A7,Credits!`$C:`$C)"
$wsSummary.Range("D7").Formula = "=B7+C7"
$wsSummary.Range("E7").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A7,'Infrastructure Costs'!`$H:`$H)"
$wsSummary.Range("F7").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A7,'Infrastructure Costs'!`$I:`$I)"
$wsSummary.Range("G7").Formula = "=D7+E7+F7"

$wsSummary.Range("A8").Value = "Facilities"
$wsSummary.Range("B8").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A8,'Infrastructure Costs'!`$G:`$G)"
$wsSummary.Range("C8").Formula = "=SUMIF(Credits!`$A:`$A,A8,Credits!`$C:`$C)"
$wsSummary.Range("D8").Formula = "=B8+C8"
$wsSummary.Range("E8").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A8,'Infrastructure Costs'!`$H:`$H)"
$wsSummary.Range("F8").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A8,'Infrastructure Costs'!`$I:`$I)"
$wsSummary.Range("G8").Formula = "=D8+E8+F8"

# Row 9: TOTAL (now shifted down from row 7)
$wsSummary.Range("B9").Formula = "=SUM(B3:B8)"
$wsSummary.Range("C9").Formula = "=SUM(C3:C8)"
$wsSummary.Range("D9").Formula = "=SUM(D3:D8)"
$wsSummary.Range("E9").Formula = "=SUM(E3:E8)"
$wsSummary.Range("F9").Formula = "=SUM(F3:F8)"
$wsSummary.Range("G9").Formula = "=SUM(G3:G8)"

$wsSummary.Range("A2:G9").AutoFilter()

# ---------------------------------------------------------------
# 6) Keep the workbook-level _FilterDatabase defined names in sync
#    with the new autoFilter ranges on each sheet.
# ---------------------------------------------------------------
$wb.Names.Item("Sizing Guidelines!_FilterDatabase").RefersTo = "='Sizing Guidelines'!`$A`$2:`$F`$11"
$wb.Names.Item("Infrastructure Costs!_FilterDatabase").RefersTo = "='Infrastructure Costs'!`$A`$2:`$K`$13"
$wb.Names.Item("Credits!_FilterDatabase").RefersTo = "='Credits'!`$A`$2:`$D`$8"
$wb.Names.Item("3-Year Summary!_FilterDatabase").RefersTo = "='3-Year Summary'!`$A`$2:`$G`$9"
